$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.818.95"
Set-TextValue $ws.Range("E2") "  -0.23%  "
Set-TextValue $ws.Range("D3") "3.793.89"
Set-TextValue $ws.Range("E3") "  -2.15%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "597.91"
Set-TextValue $ws.Range("E5") "  -0.19%  "
Set-TextValue $ws.Range("D6") "168.46"
Set-TextValue $ws.Range("E6") "  -1.69%  "
Set-TextValue $ws.Range("D7") "3.792.02"
Set-TextValue $ws.Range("E7") "  -2.14%  "
Set-TextValue $ws.Range("E8") "  +0.00%  "
Set-TextValue $ws.Range("E9") "  -0.07%  "
Set-TextValue $ws.Range("E10") "  +0.86%  "
Set-TextValue $ws.Range("E11") "  +1.53%  "
Set-TextValue $ws.Range("D12") "0.459"
Set-TextValue $ws.Range("E12") "  +0.60%  "
Set-TextValue $ws.Range("E13") "  +5.66%  "
Set-TextValue $ws.Range("D14") "36.79"
Set-TextValue $ws.Range("E14") "  -0.53%  "
Set-TextValue $ws.Range("D15") "4.432.36"
Set-TextValue $ws.Range("E15") "  -2.22%  "
Set-TextValue $ws.Range("D16") "3.815.64"
Set-TextValue $ws.Range("E16") "  -1.69%  "
Set-TextValue $ws.Range("D17") "18.90"
Set-TextValue $ws.Range("E17") "  +4.54%  "
Set-TextValue $ws.Range("D18") "67.827.08"
Set-TextValue $ws.Range("E18") "  -0.52%  "
Set-TextValue $ws.Range("D19") "7.31"
Set-TextValue $ws.Range("E19") "  -0.50%  "
Set-TextValue $ws.Range("E20") "  +0.63%  "
Set-TextValue $ws.Range("D21") "10.59"
Set-TextValue $ws.Range("E21") "  -2.30%  "
Set-TextValue $ws.Range("D22") "468.01"
Set-TextValue $ws.Range("E22") "  +0.25%  "
Set-TextValue $ws.Range("D23") "0.729"
Set-TextValue $ws.Range("E23") "  -1.22%  "
Set-TextValue $ws.Range("E24") "  -5.86%  "
Set-TextValue $ws.Range("D25") "83.42"
Set-TextValue $ws.Range("E25") "  +0.19%  "
Set-TextValue $ws.Range("D26") "2.28"
Set-TextValue $ws.Range("E26") "  +2.04%  "
Set-TextValue $ws.Range("D27") "12.19"
Set-TextValue $ws.Range("E27") "  +1.00%  "
Set-TextValue $ws.Range("D28") "10.29"
Set-TextValue $ws.Range("E28") "  +3.27%  "
Set-TextValue $ws.Range("E29") "  -0.02%  "
Set-TextValue $ws.Range("E30") "  -0.85%  "
Set-TextValue $ws.Range("D31") "3.949.11"
Set-TextValue $ws.Range("E31") "  -2.07%  "
Set-TextValue $ws.Range("D32") "7.66"
Set-TextValue $ws.Range("E32") "  -0.38%  "
Set-TextValue $ws.Range("E33") "  -2.25%  "
Set-TextValue $ws.Range("D34") "30.52"
Set-TextValue $ws.Range("E34") "  -2.26%  "
Set-TextValue $ws.Range("D35") "9.21"
Set-TextValue $ws.Range("E35") "  -2.85%  "
Set-TextValue $ws.Range("D36") "3.757.09"
Set-TextValue $ws.Range("E36") "  -2.47%  "
Set-TextValue $ws.Range("E37") "  +0.68%  "
Set-TextValue $ws.Range("E38") "  -1.44%  "
Set-TextValue $ws.Range("D39") "5.92"
Set-TextValue $ws.Range("E39") "  +0.44%  "
Set-TextValue $ws.Range("E40") "  -1.56%  "
Set-TextValue $ws.Range("E41") "  -1.40%  "
Set-TextValue $ws.Range("D42") "1.00"
Set-TextValue $ws.Range("E42") "  +0.02%  "
Set-TextValue $ws.Range("D43") "0.317"
Set-TextValue $ws.Range("E43") "  +1.18%  "
Set-TextValue $ws.Range("D45") "8.76"
Set-TextValue $ws.Range("E45") "  +1.52%  "
Set-TextValue $ws.Range("E46") "  -0.66%  "
Set-TextValue $ws.Range("D47") "409.20"
Set-TextValue $ws.Range("E47") "  -3.73%  "
Set-TextValue $ws.Range("D48") "46.28"
Set-TextValue $ws.Range("E48") "  -1.98%  "
Set-TextValue $ws.Range("E49") "  -7.84%  "
Set-TextValue $ws.Range("D50") "142.13"
Set-TextValue $ws.Range("E50") "  -0.83%  "
Set-TextValue $ws.Range("D51") "0.0355"
Set-TextValue $ws.Range("E51") "  -0.68%  "
